# Update countries & provincias Spain
# Applies the daily COVID data refresh described in the commit:
#  - updates the "Datos actualizados" timestamp
#  - updates case/death/recovered counters for a handful of countries
#  - re-ranks three pairs/triples of countries whose totals crossed each
#    other in the sort order (El Salvador/Costa Rica/Australia,
#    Hong Kong/Montenegro/Libia, Groenlandia/Islas Malvinas)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 11:57"

# --- Straight numeric refreshes (country stays on its row) -------------
# Row 6: India
$ws.Range("B6").Value = 1389221
$ws.Range("C6").Value = 3727
$ws.Range("E6").Value = 469797

# Row 20: Banglades
$ws.Range("B20").Value = 223453
$ws.Range("C20").Value = 2275
$ws.Range("D20").Value = 123882
$ws.Range("E20").Value = 96643
$ws.Range("G20").Value = 54
$ws.Range("H20").Value = 2928

# Row 27: Indonesia
$ws.Range("B27").Value = 98778
$ws.Range("C27").Value = 1492
$ws.Range("D27").Value = 56655
$ws.Range("E27").Value = 37342
$ws.Range("G27").Value = 67
$ws.Range("H27").Value = 4781

# Row 34: Oman
$ws.Range("B34").Value = 76005
$ws.Range("C34").Value = 1147
$ws.Range("D34").Value = 55299
$ws.Range("E34").Value = 20322
$ws.Range("G34").Value = 13
$ws.Range("H34").Value = 384

# Row 49: Polonia
$ws.Range("B49").Value = 43065
$ws.Range("C49").Value = 443
$ws.Range("D49").Value = 32753
$ws.Range("E49").Value = 8641
$ws.Range("G49").Value = 7
$ws.Range("H49").Value = 1671

# Row 51: Barein
$ws.Range("E51").Value = 3402
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 140

# Row 54: Afganistan
$ws.Range("B54").Value = 36157
$ws.Range("C54").Value = 121
$ws.Range("D54").Value = 25180
$ws.Range("E54").Value = 9718
$ws.Range("G54").Value = 11
$ws.Range("H54").Value = 1259

# Row 86: Malasia
$ws.Range("B86").Value = 8897
$ws.Range("C86").Value = 13
$ws.Range("D86").Value = 8600
$ws.Range("E86").Value = 173
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 124

# Row 89: Finlandia
$ws.Range("B89").Value = 7393
$ws.Range("C89").Value = 5
$ws.Range("E89").Value = 144

# Row 114: Sri Lanka
$ws.Range("D114").Value = 2106
$ws.Range("E114").Value = 653

# Row 126: Estonia
$ws.Range("B126").Value = 2034
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 1922
$ws.Range("E126").Value = 43

# Row 167: Birmania
$ws.Range("B167").Value = 349
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 55

# --- Re-ranked countries -------------------------------------------------
# El Salvador overtakes Costa Rica and Australia (rows 73-75)
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 14630
$ws.Range("C73").Value = 409
$ws.Range("D73").Value = 7648
$ws.Range("E73").Value = 6582
$ws.Range("G73").Value = 10
$ws.Range("H73").Value = 400

$ws.Range("A74").Value = "Costa Rica"
$ws.Range("B74").Value = 14600
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 3640
$ws.Range("E74").Value = 10862
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 98

$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 14403
$ws.Range("C75").Value = 453
$ws.Range("D75").Value = 9019
$ws.Range("E75").Value = 5229
$ws.Range("G75").Value = 10
$ws.Range("H75").Value = 155

# Hong Kong overtakes Montenegro and Libia (rows 115-117)
$ws.Range("A115").Value = "Hong Kong"
$ws.Range("B115").Value = 2762
$ws.Range("C115").Value = 256
$ws.Range("D115").Value = 1495
$ws.Range("E115").Value = 1249
$ws.Range("H115").Value = 18

$ws.Range("A116").Value = "Montenegro"
$ws.Range("B116").Value = 2747
$ws.Range("D116").Value = 664
$ws.Range("E116").Value = 2040
$ws.Range("H116").Value = 43

$ws.Range("A117").Value = "Libia"
$ws.Range("B117").Value = 2547
$ws.Range("D117").Value = 510
$ws.Range("E117").Value = 1979
$ws.Range("H117").Value = 58

# Groenlandia/Islas Malvinas swap places (identical totals, tie-break change)
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
